$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data's "batsman" value has a trailing NO-BREAK SPACE
# (U+00A0) rather than a plain space (matches the existing cells already
# present in the workbook byte-for-byte).
$batsman = "Jos Buttler" + [char]0x00A0

# Helper: write $val into $addr while forcing it to be stored as text
# (OOXML t="str"), even when $val looks like a number (e.g. "24",
# "96.00"). A plain "$ws.Range($addr).Value = $val" would otherwise be
# silently coerced to a numeric cell by the Value setter.
# We stage the text in an always-out-of-range scratch cell formatted as
# Text, copy it, and paste-special *values* into the destination so the
# destination cell's own style/format is left untouched.
function Set-TextValue($ws, $addr, $val) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- 1. Insert two new columns (ownTeam, oppTeam) before the existing
#        "batsman" column (D), shifting batsman..sr from D:I to F:K ---
$ws.Range("D1:E1").EntireColumn.Insert()

# --- 2. Insert a new row for the "Abu Dhabi / October 30 2020" match
#        above the current row 3 ("Dubai (DSC) / October 22 2020"),
#        shifting it down to row 4 ---
$ws.Range("A3").EntireRow.Insert()

# --- 3. Headers ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# --- 4. Row 2 (existing match vs RCB) - add team columns ---
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Royal Challengers Bangalore"

# --- 5. Row 3 (new match vs Kings XI Punjab) ---
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 30 2020"
$ws.Range("C3").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Kings XI Punjab"
$ws.Range("F3").Value = $batsman
Set-TextValue $ws "G3" "22"
Set-TextValue $ws "H3" "11"
Set-TextValue $ws "I3" "1"
Set-TextValue $ws "J3" "2"
Set-TextValue $ws "K3" "200.00"

# --- 6. Row 4 (existing match vs Sunrisers, shifted down from row 3) ---
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Sunrisers Hyderabad"

# --- 7. Row 5 (new match vs Kolkata Knight Riders) ---
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " November 01 2020"
$ws.Range("C5").Value = "KKR won by 60 runs"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Kolkata Knight Riders"
$ws.Range("F5").Value = $batsman
Set-TextValue $ws "G5" "35"
Set-TextValue $ws "H5" "22"
Set-TextValue $ws "I5" "4"
Set-TextValue $ws "J5" "1"
Set-TextValue $ws "K5" "159.09"

# --- 8. Row 6 (new match vs Chennai Super Kings) ---
$ws.Range("A6").Value = " Abu Dhabi"
$ws.Range("B6").Value = " October 19 2020"
$ws.Range("C6").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D6").Value = "Rajasthan Royals"
$ws.Range("E6").Value = "Chennai Super Kings"
$ws.Range("F6").Value = $batsman
Set-TextValue $ws "G6" "70"
Set-TextValue $ws "H6" "48"
Set-TextValue $ws "I6" "7"
Set-TextValue $ws "J6" "2"
Set-TextValue $ws "K6" "145.83"
